$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2500
$ws.Range("J2").Value = 2500
$ws.Range("L2").Value = 2500
$ws.Range("N2").Value = -2726

$ws.Range("H12").Value = 196.2
$ws.Range("I12").Value = 170.25
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 170.25
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -0.25
$ws.Range("N12").Value = -640

$ws.Range("H13").Value = 1400
$ws.Range("I13").Value = 800
$ws.Range("K13").Value = 800
$ws.Range("M13").Value = -631

$ws.Range("H33").Value = 451
$ws.Range("I33").Value = 124.23077
$ws.Range("K33").Value = 124.23077
$ws.Range("M33").Value = 104.76923

$ws.Range("H55").Value = 319.2
$ws.Range("I55").Value = 188.5
$ws.Range("J55").Value = 842
$ws.Range("K55").Value = 188.5
$ws.Range("L55").Value = 842
$ws.Range("M55").Value = 25.5
$ws.Range("N55").Value = -1270

$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws.Range("H109").Value = 85497
$ws.Range("J109").Value = 85497
$ws.Range("L109").Value = 85497
$ws.Range("N109").Value = -88271

$ws.Range("H138").Value = 3527.3845
$ws.Range("I138").Value = 2127
$ws.Range("J138").Value = 5161.1665
$ws.Range("K138").Value = 6381
$ws.Range("L138").Value = 15483.4995
$ws.Range("M138").Value = -1241
$ws.Range("N138").Value = -25763.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H16").Value = 4002.4
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 6666.6665
$ws.Range("K16").Value = 6
$ws.Range("L16").Value = 6666.6665
$ws.Range("M16").Value = 281
$ws.Range("N16").Value = -7240.6665

$ws.Range("H32").Value = 1868.5
$ws.Range("I32").Value = 1840.6
$ws.Range("J32").Value = 1952.2
$ws.Range("K32").Value = 1840.6
$ws.Range("L32").Value = 1952.2
$ws.Range("M32").Value = -1553.6
$ws.Range("N32").Value = -2526.2

$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("N88").Value = -2812

$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("N91").Value = -4808

$ws.Range("H134").Value = 61805
$ws.Range("J134").Value = 61805
$ws.Range("L134").Value = 61805
$ws.Range("N134").Value = -71945

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 9333.333000000001
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 11500
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 11500
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -12130

$ws.Range("H79").Value = 9333.333000000001
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 11500
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 11500
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -13684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258.7857
$ws.Range("I7").Value = 303.3
$ws.Range("J7").Value = 147.5
$ws.Range("K7").Value = 303.3
$ws.Range("L7").Value = 147.5
$ws.Range("M7").Value = -190.3
$ws.Range("N7").Value = -373.5

$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -14992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 300358.25
$ws.Range("I4").Value = 312875.75
$ws.Range("K4").Value = 938627.25
$ws.Range("M4").Value = -938515.25

$ws.Range("H33").Value = 20
$ws.Range("I33").Value = 20
$ws.Range("K33").Value = 120
$ws.Range("M33").Value = 163

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 333.33334
$ws.Range("J4").Value = 475
$ws.Range("L4").Value = 475
$ws.Range("N4").Value = -699

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1900
$ws.Range("I2").Value = 1900
$ws.Range("K2").Value = 1900
$ws.Range("M2").Value = -1788

$ws.Range("H40").Value = 9712.933999999999
$ws.Range("I40").Value = 9712.933999999999
$ws.Range("K40").Value = 9712.933999999999
$ws.Range("M40").Value = -9576.933999999999

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws.Range("H55").Value = 3157.3076
$ws.Range("I55").Value = 2443.5
$ws.Range("J55").Value = 4299.4
$ws.Range("K55").Value = 2443.5
$ws.Range("L55").Value = 4299.4
$ws.Range("M55").Value = -2270.5
$ws.Range("N55").Value = -4645.4

$ws.Range("H127").Value = 99994
$ws.Range("J127").Value = 99994
$ws.Range("L127").Value = 99994
$ws.Range("N127").Value = -109914

$ws.Range("H132").Value = 5386.5835
$ws.Range("I132").Value = 5454.875
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 16364.625
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -13834.625
$ws.Range("N132").Value = -20810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 113526.2
$ws.Range("J2").Value = 12631
$ws.Range("L2").Value = 12631
$ws.Range("N2").Value = -12855

$ws.Range("H63").Value = 40000
$ws.Range("I63").Value = 40000
$ws.Range("K63").Value = 40000
$ws.Range("M63").Value = -39376

$ws.Range("H66").Value = 40000
$ws.Range("I66").Value = 40000
$ws.Range("K66").Value = 120000
$ws.Range("M66").Value = -116880

$ws.Range("H113").Value = 519.375
$ws.Range("I113").Value = 429.33334
$ws.Range("J113").Value = 789.5
$ws.Range("K113").Value = 1288.00002
$ws.Range("L113").Value = 2368.5
$ws.Range("M113").Value = 881.9999800000001
$ws.Range("N113").Value = -6708.5

Write-Output "applied all Kraken_Profits updates"